$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws3 = $wb.Worksheets.Item(3)

# --- Sheet1: add header "Ranking " in D1 ---
$ws1.Range("D1").Value = "Ranking "

# --- Sheet1: add new row 8 data (8/8/2019 entry) ---
$ws1.Range("A8").Value = 43685
[void]$ws1.Range("A7").Copy()
[void]$ws1.Range("A8").PasteSpecial(-4122)
$ws1.Range("B8").Value = 132477
$ws1.Range("C8").Value = 112720
$ws1.Range("D8").Value = 126722
$ws1.Range("E8").Value = "Yes"

# --- Sheet1: page setup (portrait orientation) ---
$ws1.PageSetup.Orientation = 1

# --- Sheet3: update progress count + recalculated total ---
$ws3.Range("B2").Value = 18

# --- Update selections (Sheet1 stays the active tab) ---
[void]$ws3.Range("C19").Select()
[void]$ws1.Range("C15").Select()

Write-Host "edit applied"
